$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '34.516.13'
Set-TextValue 'E2' '  -0.17%  '
Set-TextValue 'D3' '1.807.88'
Set-TextValue 'E3' '  -0.52%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '228.03'
Set-TextValue 'E5' '  -0.13%  '
Set-TextValue 'D6' '0.579'
Set-TextValue 'E6' '  +3.47%  '
Set-TextValue 'E7' '  +0.09%  '
Set-TextValue 'D8' '36.73'
Set-TextValue 'E8' '  +5.79%  '
Set-TextValue 'E9' '  -0.64%  '
Set-TextValue 'D10' '0.0694'
Set-TextValue 'E10' '  -0.29%  '
Set-TextValue 'D11' '0.0965'
Set-TextValue 'E11' '  +1.43%  '
Set-TextValue 'D12' '2.066.72'
Set-TextValue 'E12' '  -0.57%  '
Set-TextValue 'E13' '  +1.83%  '
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.814.25'
Set-TextValue 'E14' '  -0.22%  '
Set-TextValue 'B15' 'Polygon'
Set-TextValue 'C15' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D15' '0.654'
Set-TextValue 'E15' '  +1.26%  '
Set-TextValue 'D16' '4.48'
Set-TextValue 'E16' '  +3.14%  '
Set-TextValue 'D17' '34.479.78'
Set-TextValue 'E17' '  -0.30%  '
Set-TextValue 'D18' '69.65'
Set-TextValue 'E18' '  +0.72%  '
Set-TextValue 'D19' '245.52'
Set-TextValue 'E19' '  -0.74%  '
Set-TextValue 'E20' '  -1.33%  '
Set-TextValue 'D21' '11.65'
Set-TextValue 'E21' '  +0.54%  '
Set-TextValue 'E22' '  +0.10%  '
Set-TextValue 'E23' '  -0.39%  '
Set-TextValue 'D24' '2.20'
Set-TextValue 'E24' '  +5.17%  '
Set-TextValue 'D25' '172.23'
Set-TextValue 'E25' '  -0.30%  '
Set-TextValue 'E26' '  +7.49%  '
Set-TextValue 'D27' '16.93'
Set-TextValue 'E27' '  +1.11%  '
Set-TextValue 'E28' '  +1.50%  '
Set-TextValue 'E29' '  +0.06%  '
Set-TextValue 'D30' '4.03'
Set-TextValue 'E30' '  -0.25%  '
Set-TextValue 'D31' '3.86'
Set-TextValue 'E31' '  +0.37%  '
Set-TextValue 'E32' '  -0.15%  '
Set-TextValue 'E33' '  -0.38%  '
Set-TextValue 'E34' '  -1.73%  '
Set-TextValue 'D35' '1.396.53'
Set-TextValue 'E35' '  -1.71%  '
Set-TextValue 'D36' '0.674'
Set-TextValue 'E36' '  -0.22%  '
Set-TextValue 'D37' '2.47'
Set-TextValue 'E37' '  -5.39%  '
Set-TextValue 'D38' '1.06'
Set-TextValue 'E38' '  +0.02%  '
Set-TextValue 'E39' '  +0.06%  '
Set-TextValue 'D40' '83.22'
Set-TextValue 'E40' '  -3.07%  '
Set-TextValue 'D41' '0.967'
Set-TextValue 'E41' '  +1.13%  '
Set-TextValue 'D42' '2.84'
Set-TextValue 'E42' '  -0.67%  '
Set-TextValue 'E43' '  +0.72%  '
Set-TextValue 'E44' '  +7.58%  '
Set-TextValue 'D45' '13.52'
Set-TextValue 'E45' '  -2.03%  '
Set-TextValue 'E46' '  -1.12%  '
Set-TextValue 'D47' '0.0500'
Set-TextValue 'E47' '  -5.12%  '
Set-TextValue 'D48' '1.969.33'
Set-TextValue 'E48' '  -0.56%  '
Set-TextValue 'D49' '104.59'
Set-TextValue 'E49' '  -1.11%  '
Set-TextValue 'E50' '  +0.07%  '
Set-TextValue 'E51' '  -3.07%  '
